$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: A23 = 5
$ws.Range("A23").Value = 5

# Update the view: scroll so A7 is the top-left visible cell, and select A24
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A24").Select()
